# "maj template comment à la fin"
# Move the "Comment" column (currently column J, with its sub-header rows)
# to the end of the J:L block, after DilutionFactor (K) and LivingCellCount
# (L). Equivalent to rotating the content of columns J, K, L left by one:
#   new J = old K   (DilutionFactor)
#   new K = old L   (LivingCellCount)
#   new L = old J   (Comment)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 5

# Capture the original values of J, K, L for every row first, since we will
# overwrite them in place.
$origJ = @{}
$origK = @{}
$origL = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $origJ[$r] = $ws.Cells.Item($r, 10).Value2
    $origK[$r] = $ws.Cells.Item($r, 11).Value2
    $origL[$r] = $ws.Cells.Item($r, 12).Value2
}

for ($r = 1; $r -le $lastRow; $r++) {
    $newJ = $origK[$r]
    $newK = $origL[$r]
    $newL = $origJ[$r]

    $ws.Cells.Item($r, 10).Value2 = $newJ
    $ws.Cells.Item($r, 11).Value2 = $newK
    $ws.Cells.Item($r, 12).Value2 = $newL
}
